# Update "want to go" headcount (column F) figures refreshed at commit 456a3b4
$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "展览"; Cell = "F3"; Old = 99; New = 102 }
    @{ Sheet = "展览"; Cell = "F4"; Old = 925; New = 930 }
    @{ Sheet = "展览"; Cell = "F5"; Old = 64; New = 65 }
    @{ Sheet = "展览"; Cell = "F6"; Old = 7203; New = 7231 }
    @{ Sheet = "展览"; Cell = "F8"; Old = 161; New = 162 }
    @{ Sheet = "展览"; Cell = "F9"; Old = 6579; New = 6591 }
    @{ Sheet = "展览"; Cell = "F10"; Old = 134; New = 135 }
    @{ Sheet = "展览"; Cell = "F11"; Old = 278; New = 281 }
    @{ Sheet = "展览"; Cell = "F12"; Old = 4549; New = 4568 }
    @{ Sheet = "展览"; Cell = "F16"; Old = 4628; New = 4654 }
    @{ Sheet = "展览"; Cell = "F18"; Old = 250; New = 256 }
    @{ Sheet = "展览"; Cell = "F20"; Old = 354; New = 356 }
    @{ Sheet = "展览"; Cell = "F21"; Old = 235; New = 236 }
    @{ Sheet = "展览"; Cell = "F23"; Old = 175; New = 177 }
    @{ Sheet = "展览"; Cell = "F28"; Old = 8260; New = 8275 }
    @{ Sheet = "展览"; Cell = "F30"; Old = 1434; New = 1442 }
    @{ Sheet = "展览"; Cell = "F32"; Old = 729; New = 730 }
    @{ Sheet = "展览"; Cell = "F34"; Old = 53; New = 54 }
    @{ Sheet = "展览"; Cell = "F35"; Old = 990; New = 991 }
    @{ Sheet = "展览"; Cell = "F37"; Old = 1696; New = 1705 }
    @{ Sheet = "展览"; Cell = "F39"; Old = 978; New = 982 }
    @{ Sheet = "展览"; Cell = "F40"; Old = 43; New = 44 }
    @{ Sheet = "展览"; Cell = "F41"; Old = 4291; New = 4319 }
    @{ Sheet = "展览"; Cell = "F42"; Old = 357; New = 358 }
    @{ Sheet = "展览"; Cell = "F44"; Old = 120; New = 122 }
    @{ Sheet = "展览"; Cell = "F47"; Old = 1136; New = 1139 }
    @{ Sheet = "展览"; Cell = "F49"; Old = 29; New = 31 }
    @{ Sheet = "本地生活"; Cell = "F2"; Old = 241; New = 242 }
    @{ Sheet = "全部类型"; Cell = "F2"; Old = 241; New = 242 }
    @{ Sheet = "全部类型"; Cell = "F6"; Old = 99; New = 102 }
    @{ Sheet = "全部类型"; Cell = "F8"; Old = 925; New = 930 }
    @{ Sheet = "全部类型"; Cell = "F9"; Old = 64; New = 65 }
    @{ Sheet = "全部类型"; Cell = "F10"; Old = 7203; New = 7231 }
    @{ Sheet = "全部类型"; Cell = "F12"; Old = 161; New = 162 }
    @{ Sheet = "全部类型"; Cell = "F13"; Old = 6579; New = 6591 }
    @{ Sheet = "全部类型"; Cell = "F14"; Old = 134; New = 135 }
    @{ Sheet = "全部类型"; Cell = "F15"; Old = 278; New = 281 }
    @{ Sheet = "全部类型"; Cell = "F16"; Old = 4549; New = 4568 }
    @{ Sheet = "全部类型"; Cell = "F20"; Old = 4628; New = 4654 }
    @{ Sheet = "全部类型"; Cell = "F22"; Old = 250; New = 256 }
    @{ Sheet = "全部类型"; Cell = "F23"; Old = 354; New = 356 }
    @{ Sheet = "全部类型"; Cell = "F24"; Old = 235; New = 236 }
    @{ Sheet = "全部类型"; Cell = "F29"; Old = 8260; New = 8275 }
    @{ Sheet = "全部类型"; Cell = "F31"; Old = 1434; New = 1442 }
    @{ Sheet = "全部类型"; Cell = "F33"; Old = 729; New = 730 }
    @{ Sheet = "全部类型"; Cell = "F35"; Old = 53; New = 54 }
    @{ Sheet = "全部类型"; Cell = "F36"; Old = 990; New = 991 }
    @{ Sheet = "全部类型"; Cell = "F37"; Old = 1696; New = 1705 }
    @{ Sheet = "全部类型"; Cell = "F39"; Old = 978; New = 982 }
    @{ Sheet = "全部类型"; Cell = "F40"; Old = 43; New = 44 }
    @{ Sheet = "全部类型"; Cell = "F41"; Old = 4291; New = 4319 }
    @{ Sheet = "全部类型"; Cell = "F42"; Old = 357; New = 358 }
    @{ Sheet = "全部类型"; Cell = "F44"; Old = 120; New = 122 }
    @{ Sheet = "全部类型"; Cell = "F47"; Old = 1136; New = 1139 }
    @{ Sheet = "全部类型"; Cell = "F49"; Old = 29; New = 31 }
)

$mismatchCount = 0
foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $rng = $ws.Range($u.Cell)
    if ($rng.Value2 -ne $u.Old) {
        $mismatchCount++
        Write-Output "Unexpected existing value on $($u.Sheet)!$($u.Cell): $($rng.Value2) (expected $($u.Old))"
    }
    $rng.Value = $u.New
}

Write-Output "Applied $($updates.Count) cell updates ($mismatchCount unexpected prior values)."

